$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Website URL"
$ws.Range("B1").Value = "Company name"
$ws.Range("C1").Value = "Record ID"
